$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha, serial), M (Volumen), N/O/P (Precio min/max/ponderado), S (Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44438; M = 60;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 3;  D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 4;  D = 44476; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 5;  D = 44357; M = 35;  N = 1000; O = 1000; P = 1000; S = 1000 },
    @{ Row = 6;  D = 44473; M = 120; N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 7;  D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 8;  D = 44432; M = 30;  N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 9;  D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 },
    @{ Row = 10; D = 44418; M = 40;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 11; D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 12; D = 44417; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 },
    @{ Row = 13; D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value2 = $r.D   # D = Fecha
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M = Volumen
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N = Precio minimo
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O = Precio maximo
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P = Precio promedio ponderado
    $ws.Cells.Item($r.Row, 19).Value = $r.S   # S = Precio $/Kg
}
